# summer 24 week 11 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.26
$ws.Range("E3").Value = 1.31
$ws.Range("F3").Value = 1.16
$ws.Range("B4").Value = 1.49
$ws.Range("E4").Value = 1.23
$ws.Range("C5").Value = 1.36
$ws.Range("D5").Value = 1.33
$ws.Range("C6").Value = 1.55
$ws.Range("G6").Value = 0.96
